$d = $word.ActiveDocument

# 1. "Learning Objectives" heading is split across two runs: "Lear" + "ning Objectives".
#    The edit re-splits it as "Learn" + "ing Objectives" (the combined text stays the same).
$d.Content.Find.Execute("ning Objectives", $true, $false, $false, $false, $false, $true, 1, $false, "ing Objectives", 2) | Out-Null
$d.Content.Find.Execute("Lear ", $true, $false, $false, $false, $false, $true, 1, $false, "Learn ", 2) | Out-Null

# 2. Course code header "VE 3500" -> "VE3500" (remove the stray space before "3").
foreach ($sec in $d.Sections) {
    $hdr = $sec.Headers.Item(1)
    $hdr.Range.Find.Execute(" 3500", $true, $false, $false, $false, $false, $true, 1, $false, "3500", 2) | Out-Null
}
